$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "id" header in column A and lowercase the existing headers
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "tv"
$ws.Range("C1").Value = "radio"
$ws.Range("D1").Value = "newspaper"
$ws.Range("E1").Value = "sales"
